$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'107"
$ws.Range("D2").Value = "'106"
$ws.Range("C3").Value = "'5990"
$ws.Range("D3").Value = "'5985"
$ws.Range("C8").Value = "'2.6"
$ws.Range("D8").Value = "'2.4"
$ws.Range("H8").Value = "Oct/24"
$ws.Range("H9").Value = "Oct/24"
$ws.Range("C67").Value = "'6.4"
$ws.Range("D67").Value = "'6.5"
$ws.Range("H67").Value = "Sep/24"
$ws.Range("C87").Value = "'2.6"
$ws.Range("D87").Value = "'2.4"
$ws.Range("H87").Value = "Oct/24"
$ws.Range("H88").Value = "Oct/24"
$ws.Range("C89").Value = "'316"
$ws.Range("H89").Value = "Oct/24"
$ws.Range("C90").Value = "'322"
$ws.Range("D90").Value = "'321"
$ws.Range("H90").Value = "Oct/24"
$ws.Range("D91").Value = "'3.3"
$ws.Range("H91").Value = "Oct/24"
$ws.Range("C97").Value = "'2.1"
$ws.Range("D97").Value = "'2.3"
$ws.Range("H97").Value = "Oct/24"
$ws.Range("H98").Value = "Oct/24"
$ws.Range("D104").Value = "'315"
$ws.Range("H104").Value = "Oct/24"
$ws.Range("C105").Value = "'-4.9"
$ws.Range("D105").Value = "'-6.8"
$ws.Range("H105").Value = "Oct/24"
$ws.Range("D119").Value = "'4.9"
$ws.Range("H119").Value = "Oct/24"
$ws.Range("D120").Value = "'4.7"
$ws.Range("H120").Value = "Oct/24"
$ws.Range("C123").Value = "'2.3"
$ws.Range("D123").Value = "'2.4"
$ws.Range("H123").Value = "Oct/24"
$ws.Range("D124").Value = "'337"
$ws.Range("H124").Value = "Oct/24"
$ws.Range("D126").Value = "'270"
$ws.Range("H126").Value = "Oct/24"
$ws.Range("C128").Value = "'2.9"
$ws.Range("H128").Value = "Oct/24"
$ws.Range("C139").Value = "'4.58"
$ws.Range("D139").Value = "'4.58"
$ws.Range("C141").Value = "'4.92"
$ws.Range("D141").Value = "'4.71"
$ws.Range("H141").Value = "Oct/24"
$ws.Range("C142").Value = "'4.6"
$ws.Range("D151").Value = "'108"
$ws.Range("H151").Value = "Sep/24"
$ws.Range("C162").Value = "'8391"
$ws.Range("D162").Value = "'9681"
$ws.Range("H162").Value = "Sep/24"
$ws.Range("C169").Value = "'-257000"
$ws.Range("D169").Value = "'64000"
$ws.Range("H169").Value = "Oct/24"
$ws.Range("C198").Value = "'22762"
$ws.Range("D198").Value = "'22060"
$ws.Range("H198").Value = "Sep/24"
$ws.Range("C200").Value = "'93.7"
$ws.Range("D200").Value = "'91.5"
$ws.Range("H200").Value = "Oct/24"
$ws.Range("C224").Value = "'42.3"
$ws.Range("D224").Value = "'44.1"
$ws.Range("H224").Value = "Oct/24"
$ws.Range("C226").Value = "'42.6"
$ws.Range("D226").Value = "'43.9"
$ws.Range("H226").Value = "Oct/24"
$ws.Range("C229").Value = "'46.2"
$ws.Range("D229").Value = "'49.8"
$ws.Range("H229").Value = "Oct/24"
$ws.Range("C230").Value = "'52"
$ws.Range("D230").Value = "'52.2"
$ws.Range("H230").Value = "Oct/24"
$ws.Range("C275").Value = "'12556"
$ws.Range("D275").Value = "'12485"
$ws.Range("H275").Value = "Oct/24"
$ws.Range("D276").Value = "'8"
$ws.Range("H276").Value = "Oct/24"
$ws.Range("C277").Value = "'53.2"
$ws.Range("D277").Value = "'46.9"
$ws.Range("H277").Value = "Nov/24"
$ws.Range("C278").Value = "'4.8"
$ws.Range("D278").Value = "'6"
$ws.Range("C279").Value = "'601"
$ws.Range("D279").Value = "'599"
$ws.Range("H279").Value = "Sep/24"
$ws.Range("C280").Value = "'1.64"
$ws.Range("D280").Value = "'1.63"
$ws.Range("H280").Value = "Sep/24"
$ws.Range("C281").Value = "'1.17"
$ws.Range("D281").Value = "'1.14"
$ws.Range("H281").Value = "Sep/24"
$ws.Range("C282").Value = "'12.59"
$ws.Range("D282").Value = "'12.52"
$ws.Range("H282").Value = "Sep/24"
$ws.Range("C283").Value = "'1.61"
$ws.Range("D283").Value = "'1.59"
$ws.Range("H283").Value = "Sep/24"
$ws.Range("C284").Value = "'17.9"
$ws.Range("D284").Value = "'17.8"
$ws.Range("H284").Value = "Sep/24"
$ws.Range("C301").Value = "'6.86"
$ws.Range("D301").Value = "'6.81"
$ws.Range("C302").Value = "'0.5"
$ws.Range("D302").Value = "'-10.8"
$ws.Range("C316").Value = "'192"
$ws.Range("D316").Value = "'191"
$ws.Range("C317").Value = "'506"
$ws.Range("D317").Value = "'514"
$ws.Range("C318").Value = "'133"
$ws.Range("D318").Value = "'131"
$ws.Range("C319").Value = "'448"
$ws.Range("D319").Value = "'374"
$ws.Range("H319").Value = "Sep/24"
$ws.Range("C332").Value = "'-0.78"
$ws.Range("D332").Value = "'3.13"
$ws.Range("C338").Value = "'-1.86"
$ws.Range("D338").Value = "'1.72"
$ws.Range("C339").Value = "'1.14"
$ws.Range("D339").Value = "'-0.85"
$ws.Range("C340").Value = "'0.31"
$ws.Range("D340").Value = "'-0.93"
